$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "data/damages_EAD/static/hazard/RP_100.tif"
$ws.Range("D3").Value = "data/damages_EAD/static/hazard/RP_100.tif"
$ws.Range("D4").Value = "data/damages_EAD/static/hazard/RP_1000.tif"
$ws.Range("D5").Value = "data/damages_EAD/static/hazard/RP_1000.tif"
